$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("ViolMetTableKTR")
$ws3 = $wb.Worksheets.Item("ViolMetTableNFkB")
$ws4 = $wb.Worksheets.Item("PeakStats")
$ws5 = $wb.Worksheets.Item("SignalStats")

# --- Step 1: rename the existing A12 labels back to their base feature names ---
# (this frees up the old unique shared strings "duration_ktr(:,3)" / "duration_nfkb(:,3)")
$ws2.Range("A12").Value = "duration_ktr"
$ws3.Range("A12").Value = "duration_nfkb"

# --- Step 2: add the shared "column Index" header in column B (row 1) on both sheets ---
$ws2.Range("B1").Value = "column Index"
$ws3.Range("B1").Value = "column Index"

# --- Step 3: fill in column B counts for the existing rows (2-11) ---
for ($r = 2; $r -le 11; $r++) {
    $ws2.Cells.Item($r, 2).Value = 1
    $ws3.Cells.Item($r, 2).Value = 1
}

# row 12 (duration_ktr / duration_nfkb) gets a count of 12
$ws2.Cells.Item(12, 2).Value = 12
$ws3.Cells.Item(12, 2).Value = 12

# --- Step 4: append the new feature rows (13-18) on ViolMetTableKTR ---
$ws2.Range("A13").Value = "max_peak2trough_ktr"
$ws2.Range("B13").Value = 1

$ws2.Range("A14").Value = "pk2_prom_ktr"
$ws2.Range("B14").Value = 1

$ws2.Range("A15").Value = "pk1_prom_ktr "
$ws2.Range("B15").Value = 1

$ws2.Range("A16").Value = "pk2_width_ktr "
$ws2.Range("B16").Value = 1

$ws2.Range("A17").Value = "pk1_width_ktr "
$ws2.Range("B17").Value = 1

$ws2.Range("A18").Value = "medfreq_ktr      "
$ws2.Range("B18").Value = 1

# --- Step 5: append the new feature rows (13-18) on ViolMetTableNFkB ---
$ws3.Range("A13").Value = "max_peak2trough_nfkb     "
$ws3.Range("B13").Value = 1

$ws3.Range("A14").Value = "pk2_prom_nfkb"
$ws3.Range("B14").Value = 1

$ws3.Range("A15").Value = "pk1_prom_nfkb "
$ws3.Range("B15").Value = 1

$ws3.Range("A16").Value = "pk2_width_nfkb "
$ws3.Range("B16").Value = 1

$ws3.Range("A17").Value = "pk1_width_nfkb "
$ws3.Range("B17").Value = 1

$ws3.Range("A18").Value = "medfreq_nfkb      "
$ws3.Range("B18").Value = 1

# --- Step 6: update selections to match the final interactive state ---
# (order matters: the sheet selected/activated last becomes the workbook's
# active/visible tab, so re-activate ViolMetTableNFkB at the very end to
# match the saved file's original active tab)
$ws2.Range("B14:B18").Select()
$ws4.Activate()
$ws4.Range("A36:A39").Select()
$ws5.Activate()
$ws5.Range("A17").Select()
$ws3.Activate()
$ws3.Range("C26").Select()
